$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.843.67"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.740.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.04"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5148"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2808"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.05"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06093"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.741.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06975"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.22"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6342"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.496"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.46"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.0000"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.863.76"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006583"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.961.57"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.067"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.428"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.099"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.55"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.507"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.819"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.95"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.54"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08265"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.611"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.416"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04400"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.617"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9698"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5985"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.669"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9990"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.896"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.63"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3824"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7250"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.884"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05461"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.255"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1101"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.01"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.458"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.40%  "
